$wb = $excel.ActiveWorkbook

# ---- Sheet "Results" : data changes ----
$results = $wb.Worksheets.Item("Results")

# Header block (row 1)
$results.Range("D1").Value = "4.15.0 - M0 - V2.uat"
$results.Range("G1").Value = "25/09/2024"

# Column header row (row 3) - "Script Num" column retitled "Test Data"
$results.Range("A3").Value = "Test Data"

# Old sample/demo data row (row 4) cleared out
$results.Range("A4").Value = $null
$results.Range("B4").Value = $null
$results.Range("C4").Value = $null
$results.Range("D4").Value = $null
$results.Range("E4").Value = $null

# New test-result row (row 23)
$results.Range("A23").Value = "20"
$results.Range("B23").Value = "Money express transfer"
$results.Range("C23").Value = "2024-10-03"
$results.Range("D23").Value = "EN"
$results.Range("E23").Value = " FAIL "

# ---- View changes ----
# "HTML Rpoert" sheet view (zoom), keep its own selection as-is
$htmlReport = $wb.Worksheets.Item("HTML Rpoert")
$htmlReport.Select()
$htmlWin = $excel.ActiveWindow
$htmlWin.Zoom = 130
$htmlReport.Range("F57").Select()

# Back to "Results" sheet view (zoom + new selection), restoring it as the active/visible tab
$results.Select()
$resultsWin = $excel.ActiveWindow
$resultsWin.Zoom = 40
$results.Range("E25").Select()
